$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: update the CODE action text (dummy -> variables)
$ws.Range("E2").Value = "IF [item] THEN GO(variables) ELSE (FOR item IN [cart] DO (INJECT_WITH(Explanation, item)), FINISH())"

# Row 3: rename "dummy" entry to "variables", change type from Q to CODE,
# and replace the question text with the new CODE action text (moves from D3 to E3)
$ws.Range("A3").Value = "variables"
$ws.Range("B3").Value = "Variables"
$ws.Range("C3").Value = "CODE"
$ws.Range("D3").Value = ""
$ws.Range("E3").Value = "item_api_category = TO_TEXT(SELECT a.CategoryChoice FROM Products p JOIN API a ON p.APIID == a.ID WHERE p.ID == [item]), item_med_form = TO_TEXT(SELECT MedFormID FROM Products WHERE ID == [item]), SAVE(item_api_category), SAVE(item_med_form), GO(product)"

# Row 4: add new "product" entry (type P) with a templated answer and a query
$ws.Range("A4").Value = "product"
$ws.Range("B4").Value = "Product"
$ws.Range("C4").Value = "P"
$ws.Range("D4").Value = "{SELECT [item_api_category] FROM Explanation WHERE MedFormID == [item_med_form]}"
$ws.Range("E4").Value = ""
$ws.Range("F4").Value = "SELECT * FROM Products WHERE ID == [item]"

# Row 5: update the answer text and the follow-up action
$ws.Range("C5").Value = "A"
$ws.Range("D5").Value = "Very nice :D"
$ws.Range("E5").Value = "DELETE(item_api_category), DELETE(item_med_form), FINISH()"

# Widen column D to fit the new, longer template text
$ws.Columns.Item(4).ColumnWidth = 33.8333333333333

# Row 5 grew very slightly taller after the text edit
$ws.Rows.Item(5).RowHeight = 13.85

# Move the selection to the newly edited cell
$ws.Application.Goto($ws.Range("E5"), $true)
